$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Normalize a few client names to uppercase (as they should have been) ---
$ws.Range("B8").Value  = "CARNES JOHANA"
$ws.Range("B35").Value = "LA PAMPA"
$ws.Range("B36").Value = "SAMY 2"

# --- Sort the existing data range (A2:E36) alphabetically by Cliente (column B) ---
$dataRange = $ws.Range("A2:E36")
$sortKey   = $ws.Range("B2:B36")
$dataRange.Sort($sortKey, 1, $null, $null, 1, $null, 1, 0)

# --- Renumber the Consecutivo column sequentially again (1..35) ---
for ($i = 0; $i -lt 35; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $i + 1
}

# --- Append the new debtor row at the bottom of the table ---
$ws.Range("A37").Value = 36
$ws.Range("B37").Value = "ABC"
$ws.Range("C37").Value = 46006
$ws.Range("D37").Value = 123456
$ws.Range("E37").Value = $false
